$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.474.12"
$ws.Range("E2").Value = "  +0.17%  "

# Row 3
$ws.Range("D3").Value = "1.932.25"
$ws.Range("E3").Value = "  +4.35%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.26%  "

# Row 6
$ws.Range("E6").Value = "  +0.00%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4763"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.45%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2877"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.39%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06653"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.88%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.94%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "107.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +26.91%  "

# Row 12
$ws.Range("D12").Value = "1.919.36"
$ws.Range("E12").Value = "  +3.76%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07616"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.94%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.170"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.83%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6637"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.40%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "306.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +24.58%  "

# Row 17
$ws.Range("D17").Value = "30.491.02"
$ws.Range("E17").Value = "  +0.35%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.52%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007609"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.60%  "

# Row 20
$ws.Range("E20").Value = "  -0.07%  "

# Row 21
$ws.Range("D21").Value = "2.174.96"
$ws.Range("E21").Value = "  +2.97%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.294"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.34%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.305"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.92%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.327"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.27%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.41%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.20%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.047"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.30%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1105"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.93%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.372"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.85%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.098"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.45%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.940"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.84%  "

# Row 33
$ws.Range("E33").Value = "  +4.11%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7443"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.77%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.155"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.47%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.756"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.02%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01968"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.96%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.691"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.35%  "

# Row 39
$ws.Range("E39").Value = "  +3.06%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8819"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.57%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "107.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.47%  "

# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.43%  "

# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.795"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.35%  "

# Row 44
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4208"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.56%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.287"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.53%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.217"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.03%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1216"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.29%  "

# Row 49
$ws.Range("E49").Value = "  +2.52%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05631"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.45%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3863"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.92%  "
